$d = $word.ActiveDocument

# wdLineStyleSingle = 1, wdLineWidth025pt = 2, wdColorAutomatic/Black = 0
$wdLineStyleSingle = 1
$wdLineWidth025pt = 2
$wdColorBlack = 0

foreach ($t in $d.Tables) {
    $b = $t.Borders
    $b.Enable = $true

    $b.DistanceFromTop = 0
    $b.DistanceFromBottom = 0
    $b.DistanceFromLeft = 0
    $b.DistanceFromRight = 0

    $b.OutsideLineStyle = $wdLineStyleSingle
    $b.OutsideLineWidth = $wdLineWidth025pt
    $b.OutsideColor = $wdColorBlack

    $b.InsideLineStyle = $wdLineStyleSingle
    $b.InsideLineWidth = $wdLineWidth025pt
    $b.InsideColor = $wdColorBlack
}

Write-Output "Applied borders to $($d.Tables.Count) tables"
